# Update odds values in Sheet1 to match the FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("V5").Value = 1.69
$ws.Range("V6").Value = 1.69
$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.29
$ws.Range("M9").Value = 1.05
$ws.Range("O9").Value = 1.25
$ws.Range("U15").Value = 1.5
$ws.Range("U16").Value = 1.53
$ws.Range("V16").Value = 2.38
$ws.Range("U17").Value = 1.91
$ws.Range("V17").Value = 1.8
$ws.Range("M18").Value = 1.02
$ws.Range("O18").Value = 1.13
$ws.Range("S18").Value = 1.22
$ws.Range("M19").Value = 1.04
$ws.Range("O19").Value = 1.2
$ws.Range("S19").Value = 1.3
$ws.Range("M20").Value = 1.04
$ws.Range("O20").Value = 1.22
$ws.Range("S20").Value = 1.3
$ws.Range("M21").Value = 1.02
$ws.Range("O21").Value = 1.14
$ws.Range("S21").Value = 1.22
$ws.Range("O22").Value = 1.1
$ws.Range("S22").Value = 1.17
$ws.Range("S23").Value = 1.47
$ws.Range("U23").Value = 1.91
$ws.Range("V23").Value = 1.8
$ws.Range("S25").Value = 1.37
$ws.Range("U25").Value = 1.8
$ws.Range("V25").Value = 1.95
$ws.Range("U26").Value = 1.62
$ws.Range("V27").Value = 1.73
$ws.Range("Q28").Value = 1.87
$ws.Range("R28").Value = 1.87
$ws.Range("U28").Value = 1.73
$ws.Range("G30").Value = 1.91
$ws.Range("Q30").Value = 1.94
$ws.Range("R30").Value = 1.79
$ws.Range("U30").Value = 1.8
$ws.Range("V30").Value = 1.91
$ws.Range("G31").Value = 2.32
$ws.Range("I31").Value = 2.65
$ws.Range("Q31").Value = 1.54
$ws.Range("U31").Value = 1.5
$ws.Range("I32").Value = 2.15
$ws.Range("R32").Value = 1.58
$ws.Range("U32").Value = 1.95
$ws.Range("V32").Value = 1.8
$ws.Range("K33").Value = 2.38
$ws.Range("U33").Value = 1.91
$ws.Range("V33").Value = 1.91
$ws.Range("U34").Value = 1.91
$ws.Range("V34").Value = 1.8
$ws.Range("J35").Value = 2.88
$ws.Range("Q35").Value = 1.92
$ws.Range("R35").Value = 1.82
$ws.Range("R37").Value = 1.54
